# "update 3-tier arch diagram and efforts"
# Fill in Sara's effort-tracking table (rows 20-22) which was previously
# left blank, mirroring the structure already used for Matteo's and
# Andrea's tables above it. The "Total effort" formula in C23 already
# exists (=SUM(C20:C22)) and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = [datetime]"2019-11-17"
$ws.Range("B20").Value = "Introduction"
$ws.Range("C20").Value = 1

$ws.Range("A21").Value = [datetime]"2019-11-19"
$ws.Range("B21").Value = "Overview + Components"
$ws.Range("C21").Value = 4

$ws.Range("A22").Value = [datetime]"2019-11-20"
$ws.Range("B22").Value = "Overview"
$ws.Range("C22").Value = 0.5

# Reflect where the author was looking/selecting when they saved.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E22").Select()
